# kategoria_elemzes.xlsx -- start coding the BOSCH material, add some text
#
# Adds a "Kiadvány kelte" (publication date) column on the left of the
# existing category table, fills in a grid of per-category counts for a
# couple of publication dates, and extends the date list down to row 26
# (placeholders for future publications). Also tidies up the sheet view,
# column widths, conditional formatting and page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row 1 / row 2 text is unchanged in content terms (same
#    section headers), but a new first column "Kiadvány kelte" is
#    inserted before them, so everything keeps its existing text.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Kiadvány kelte"

# ---------------------------------------------------------------------
# 2. Clear out the old two data rows (they held "1. kiadvány" /
#    "2. kiadvány" placeholders in column A) -- that data is replaced by
#    the real date-keyed grid below.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = ""
$ws.Range("A4").Value = ""

# ---------------------------------------------------------------------
# 3. Publication dates, column A, rows 3-26.
# ---------------------------------------------------------------------
$dates = @{
    3  = 41391
    4  = 41407
    5  = 41435
    6  = 41486
    7  = 41515
    8  = 41536
    9  = 41561
    10 = 41589
    11 = 41618
    12 = 41668
    13 = 41698
    14 = 41710
    15 = 41759
    16 = 41771
    17 = 41801
    18 = 41834
    19 = 41864
    20 = 41905
    21 = 41932
    22 = 41953
    23 = 42009
    24 = 42033
    25 = 42059
    26 = 42067
}
foreach ($r in $dates.Keys) {
    $ws.Range("A$r").Value = $dates[$r]
}
$ws.Range("A3:A26").NumberFormat = "m/d/yyyy"
$ws.Range("A3:A26").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Grid of per-category counts for the first five dates.
# ---------------------------------------------------------------------
$grid = @{
    "F3" = 3;  "G3" = 3;  "I3" = 3;  "J3" = 4;  "M3" = 4;  "Q3" = 3
    "E4" = 3;  "F4" = 3;  "M4" = 1;  "N4" = 1;  "Q4" = 6
    "C5" = 1;  "F5" = 1;  "M5" = 2;  "O5" = 1;  "Q5" = 1
    "C6" = 2;  "E6" = 5;  "F6" = 1;  "I6" = 6;  "N6" = 2;  "P6" = 1
    "E7" = 3;  "F7" = 2;  "G7" = 2;  "I7" = 3;  "Q7" = 3
}
foreach ($addr in $grid.Keys) {
    $ws.Range($addr).Value = $grid[$addr]
}

# ---------------------------------------------------------------------
# 5. Merge the new column A header cell across rows 1-2, centred both
#    ways (matches the rest of the header row).
# ---------------------------------------------------------------------
$ws.Range("A1:A2").HorizontalAlignment = -4108
$ws.Range("A1:A2").VerticalAlignment = -4108
$ws.Range("A1:A2").Merge()

# ---------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------
$ws.Range("A1").EntireColumn.ColumnWidth = 16.7109375
$ws.Range("N1").EntireColumn.ColumnWidth = 20.42578125
$ws.Range("P1").EntireColumn.ColumnWidth = 7.7109375
$ws.Range("Q1").EntireColumn.ColumnWidth = 8.140625

# ---------------------------------------------------------------------
# 7. Sheet view: zoom to 85%, move selection.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("G14").Select()

# ---------------------------------------------------------------------
# 8. Conditional formatting over the data grid: a 3-colour scale and a
#    2-colour scale (stacked, matching the authored workbook).
# ---------------------------------------------------------------------
$range = $ws.Range("B3:Q26")
$cf3 = $range.FormatConditions.AddColorScale(3)
$cf3.ColorScaleCriteria.Item(1).Type = -4135
$cf3.ColorScaleCriteria.Item(1).FormatColor.Color = 8109179
$cf3.ColorScaleCriteria.Item(2).Type = 3
$cf3.ColorScaleCriteria.Item(2).Value = 50
$cf3.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf3.ColorScaleCriteria.Item(3).Type = -4136
$cf3.ColorScaleCriteria.Item(3).FormatColor.Color = 7039248

$cf2 = $range.FormatConditions.AddColorScale(2)
$cf2.ColorScaleCriteria.Item(1).Type = -4135
$cf2.ColorScaleCriteria.Item(1).FormatColor.Color = 10284188
$cf2.ColorScaleCriteria.Item(2).Type = -4136
$cf2.ColorScaleCriteria.Item(2).FormatColor.Color = 8109179

# ---------------------------------------------------------------------
# 9. Page setup for printing.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
